# Populate the "Hasil" worksheet with scraped KPU (Indonesian election)
# recap data. Column D is written progressively (one write per TPS as it
# is scraped) before settling on its final, cumulative value - mirroring
# the original data-collection script's behaviour.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hasil")
$ws.Activate()

# --- Row 2: KEUDE BAKONGAN (3 TPS) ---
$ws.Range("A2").Value = "KEUDE BAKONGAN"
$ws.Range("B2").Value = 1101012001
$ws.Range("C2").Value = 873
$ws.Range("D2").Value = "[map[TPS1:[186 44 8]]]"
$ws.Range("E2").Value = 284
$ws.Range("D2").Value = "[map[TPS1:[186 44 8]] map[TPS2:[209 37 6]]]"
$ws.Range("F2").Value = 296
$ws.Range("D2").Value = "[map[TPS1:[186 44 8]] map[TPS2:[209 37 6]] map[TPS3:[202 38 7]]]"
$ws.Range("G2").Value = 293

# --- Row 5: GAMPONG DRIEN (2 TPS) ---
$ws.Range("A5").Value = "GAMPONG DRIEN"
$ws.Range("B5").Value = 1101012004
$ws.Range("C5").Value = 363
$ws.Range("D5").Value = "[map[TPS1:[144 19 2]]]"
$ws.Range("E5").Value = 182
$ws.Range("D5").Value = "[map[TPS1:[144 19 2]] map[TPS2:[153 14 4]]]"
$ws.Range("F5").Value = 181

# --- Row 6: DARUL IKHSAN (3 TPS) ---
$ws.Range("A6").Value = "DARUL IKHSAN"
$ws.Range("B6").Value = 1101012015
$ws.Range("C6").Value = 803
$ws.Range("D6").Value = "[map[TPS1:[210 33 4]]]"
$ws.Range("E6").Value = 271
$ws.Range("D6").Value = "[map[TPS1:[210 33 4]] map[TPS2:[207 31 1]]]"
$ws.Range("F6").Value = 265
$ws.Range("D6").Value = "[map[TPS1:[210 33 4]] map[TPS2:[207 31 1]] map[TPS3:[200 33 0]]]"
$ws.Range("G6").Value = 267

# --- Row 7: PADANG BEURAHAN (2 TPS) ---
$ws.Range("A7").Value = "PADANG BEURAHAN"
$ws.Range("B7").Value = 1101012016
$ws.Range("C7").Value = 549
$ws.Range("D7").Value = "[map[TPS1:[204 39 2]]]"
$ws.Range("E7").Value = 278
$ws.Range("D7").Value = "[map[TPS1:[204 39 2]] map[TPS2:[203 36 1]]]"
$ws.Range("F7").Value = 271

# --- Row 8: GAMPONG BARO (1 TPS) ---
$ws.Range("A8").Value = "GAMPONG BARO"
$ws.Range("B8").Value = 1101012017
$ws.Range("C8").Value = 260
$ws.Range("D8").Value = "[map[TPS1:[199 41 2]]]"
$ws.Range("E8").Value = 260

# --- Row 9: FAJAR HARAPAN (2 TPS) ---
$ws.Range("A9").Value = "FAJAR HARAPAN"
$ws.Range("B9").Value = 1101022001
$ws.Range("C9").Value = 517
$ws.Range("D9").Value = "[map[TPS1:[197 31 2]]]"
$ws.Range("E9").Value = 254
$ws.Range("D9").Value = "[map[TPS1:[197 31 2]] map[TPS2:[196 26 1]]]"
$ws.Range("F9").Value = 263

# --- Row 10: KRUENG BATEE (4 TPS) ---
$ws.Range("A10").Value = "KRUENG BATEE"
$ws.Range("B10").Value = 1101022002
$ws.Range("C10").Value = 1102
$ws.Range("D10").Value = "[map[TPS1:[212 16 2]]]"
$ws.Range("E10").Value = 271
$ws.Range("D10").Value = "[map[TPS1:[212 16 2]] map[TPS2:[225 12 2]]]"
$ws.Range("F10").Value = 277
$ws.Range("D10").Value = "[map[TPS1:[212 16 2]] map[TPS2:[225 12 2]] map[TPS3:[229 10 1]]]"
$ws.Range("G10").Value = 266
$ws.Range("D10").Value = "[map[TPS1:[212 16 2]] map[TPS2:[225 12 2]] map[TPS3:[229 10 1]] map[TPS4:[240 17 2]]]"
$ws.Range("H10").Value = 288
